# Shift all forecast rows forward by 15 days (adds GESS-period data: 18.09-25.09.2025 -> 03.10-10.10.2025)
# and refresh Prediction (col C) + Lookup (col D) values to match the new forecast window.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new date serial (col A), new Prediction value (col C, $null = unchanged), new Lookup text (col D)
$data = @(
    @(2, 45933, $null, "03.10.202511"),
    @(3, 45933, $null, "03.10.202512"),
    @(4, 45933, 0.199, "03.10.202513"),
    @(5, 45933, 0.271, "03.10.202514"),
    @(6, 45933, 0.271, "03.10.202515"),
    @(7, 45933, 0.219, "03.10.202516"),
    @(8, 45933, 0.155, "03.10.202517"),
    @(9, 45933, 0.055, "03.10.202518"),
    @(10, 45933, 0.012, "03.10.202519"),
    @(11, 45933, $null, "03.10.202520"),
    @(12, 45933, $null, "03.10.202521"),
    @(13, 45933, $null, "03.10.202522"),
    @(14, 45933, $null, "03.10.202523"),
    @(15, 45933, $null, "03.10.202524"),
    @(16, 45934, $null, "04.10.20251"),
    @(17, 45934, $null, "04.10.20252"),
    @(18, 45934, $null, "04.10.20253"),
    @(19, 45934, $null, "04.10.20254"),
    @(20, 45934, $null, "04.10.20255"),
    @(21, 45934, $null, "04.10.20256"),
    @(22, 45934, $null, "04.10.20257"),
    @(23, 45934, 0.037, "04.10.20258"),
    @(24, 45934, 0.269, "04.10.20259"),
    @(25, 45934, 0.558, "04.10.202510"),
    @(26, 45934, 0.633, "04.10.202511"),
    @(27, 45934, 0.926, "04.10.202512"),
    @(28, 45934, 1.197, "04.10.202513"),
    @(29, 45934, 1.38, "04.10.202514"),
    @(30, 45934, 1.28, "04.10.202515"),
    @(31, 45934, 1.024, "04.10.202516"),
    @(32, 45934, 0.7, "04.10.202517"),
    @(33, 45934, 0.214, "04.10.202518"),
    @(34, 45934, 0.018, "04.10.202519"),
    @(35, 45934, 0, "04.10.202520"),
    @(36, 45934, $null, "04.10.202521"),
    @(37, 45934, $null, "04.10.202522"),
    @(38, 45934, $null, "04.10.202523"),
    @(39, 45934, $null, "04.10.202524"),
    @(40, 45935, $null, "05.10.20251"),
    @(41, 45935, $null, "05.10.20252"),
    @(42, 45935, $null, "05.10.20253"),
    @(43, 45935, $null, "05.10.20254"),
    @(44, 45935, $null, "05.10.20255"),
    @(45, 45935, $null, "05.10.20256"),
    @(46, 45935, $null, "05.10.20257"),
    @(47, 45935, 0.222, "05.10.20258"),
    @(48, 45935, 1.594, "05.10.20259"),
    @(49, 45935, 2.721, "05.10.202510"),
    @(50, 45935, 3.435, "05.10.202511"),
    @(51, 45935, 3.945, "05.10.202512"),
    @(52, 45935, 4.091, "05.10.202513"),
    @(53, 45935, 4.125, "05.10.202514"),
    @(54, 45935, 3.688, "05.10.202515"),
    @(55, 45935, 2.788, "05.10.202516"),
    @(56, 45935, 1.728, "05.10.202517"),
    @(57, 45935, 0.483, "05.10.202518"),
    @(58, 45935, 0.052, "05.10.202519"),
    @(59, 45935, $null, "05.10.202520"),
    @(60, 45935, $null, "05.10.202521"),
    @(61, 45935, $null, "05.10.202522"),
    @(62, 45935, $null, "05.10.202523"),
    @(63, 45935, $null, "05.10.202524"),
    @(64, 45936, $null, "06.10.20251"),
    @(65, 45936, $null, "06.10.20252"),
    @(66, 45936, $null, "06.10.20253"),
    @(67, 45936, $null, "06.10.20254"),
    @(68, 45936, $null, "06.10.20255"),
    @(69, 45936, $null, "06.10.20256"),
    @(70, 45936, $null, "06.10.20257"),
    @(71, 45936, 0.063, "06.10.20258"),
    @(72, 45936, 0.575, "06.10.20259"),
    @(73, 45936, 1.101, "06.10.202510"),
    @(74, 45936, 1.55, "06.10.202511"),
    @(75, 45936, 2.006, "06.10.202512"),
    @(76, 45936, 2.079, "06.10.202513"),
    @(77, 45936, 1.39, "06.10.202514"),
    @(78, 45936, 1.278, "06.10.202515"),
    @(79, 45936, 0.965, "06.10.202516"),
    @(80, 45936, 0.532, "06.10.202517"),
    @(81, 45936, 0.141, "06.10.202518"),
    @(82, 45936, 0.013, "06.10.202519"),
    @(83, 45936, $null, "06.10.202520"),
    @(84, 45936, $null, "06.10.202521"),
    @(85, 45936, $null, "06.10.202522"),
    @(86, 45936, $null, "06.10.202523"),
    @(87, 45936, $null, "06.10.202524"),
    @(88, 45937, $null, "07.10.20251"),
    @(89, 45937, $null, "07.10.20252"),
    @(90, 45937, $null, "07.10.20253"),
    @(91, 45937, $null, "07.10.20254"),
    @(92, 45937, $null, "07.10.20255"),
    @(93, 45937, $null, "07.10.20256"),
    @(94, 45937, $null, "07.10.20257"),
    @(95, 45937, 0.066, "07.10.20258"),
    @(96, 45937, 0.631, "07.10.20259"),
    @(97, 45937, 1.485, "07.10.202510"),
    @(98, 45937, 2.29, "07.10.202511"),
    @(99, 45937, 2.722, "07.10.202512"),
    @(100, 45937, 2.543, "07.10.202513"),
    @(101, 45937, 2.329, "07.10.202514"),
    @(102, 45937, 2.016, "07.10.202515"),
    @(103, 45937, 0.947, "07.10.202516"),
    @(104, 45937, 0.562, "07.10.202517"),
    @(105, 45937, 0.137, "07.10.202518"),
    @(106, 45937, 0.019, "07.10.202519"),
    @(107, 45937, $null, "07.10.202520"),
    @(108, 45937, $null, "07.10.202521"),
    @(109, 45937, $null, "07.10.202522"),
    @(110, 45937, $null, "07.10.202523"),
    @(111, 45937, $null, "07.10.202524"),
    @(112, 45938, $null, "08.10.20251"),
    @(113, 45938, $null, "08.10.20252"),
    @(114, 45938, $null, "08.10.20253"),
    @(115, 45938, $null, "08.10.20254"),
    @(116, 45938, $null, "08.10.20255"),
    @(117, 45938, $null, "08.10.20256"),
    @(118, 45938, $null, "08.10.20257"),
    @(119, 45938, 0.064, "08.10.20258"),
    @(120, 45938, 0.516, "08.10.20259"),
    @(121, 45938, 0.975, "08.10.202510"),
    @(122, 45938, 1.939, "08.10.202511"),
    @(123, 45938, 2.644, "08.10.202512"),
    @(124, 45938, 2.555, "08.10.202513"),
    @(125, 45938, 2.123, "08.10.202514"),
    @(126, 45938, 1.312, "08.10.202515"),
    @(127, 45938, 0.792, "08.10.202516"),
    @(128, 45938, 0.31, "08.10.202517"),
    @(129, 45938, 0.097, "08.10.202518"),
    @(130, 45938, 0.012, "08.10.202519"),
    @(131, 45938, $null, "08.10.202520"),
    @(132, 45938, $null, "08.10.202521"),
    @(133, 45938, $null, "08.10.202522"),
    @(134, 45938, $null, "08.10.202523"),
    @(135, 45938, $null, "08.10.202524"),
    @(136, 45939, $null, "09.10.20251"),
    @(137, 45939, $null, "09.10.20252"),
    @(138, 45939, $null, "09.10.20253"),
    @(139, 45939, $null, "09.10.20254"),
    @(140, 45939, $null, "09.10.20255"),
    @(141, 45939, $null, "09.10.20256"),
    @(142, 45939, $null, "09.10.20257"),
    @(143, 45939, 0.083, "09.10.20258"),
    @(144, 45939, 0.584, "09.10.20259"),
    @(145, 45939, 1.592, "09.10.202510"),
    @(146, 45939, 2.321, "09.10.202511"),
    @(147, 45939, 2.784, "09.10.202512"),
    @(148, 45939, 2.957, "09.10.202513"),
    @(149, 45939, 3.09, "09.10.202514"),
    @(150, 45939, 2.702, "09.10.202515"),
    @(151, 45939, 1.992, "09.10.202516"),
    @(152, 45939, 1.008, "09.10.202517"),
    @(153, 45939, 0.233, "09.10.202518"),
    @(154, 45939, 0.021, "09.10.202519"),
    @(155, 45939, $null, "09.10.202520"),
    @(156, 45939, $null, "09.10.202521"),
    @(157, 45939, $null, "09.10.202522"),
    @(158, 45939, $null, "09.10.202523"),
    @(159, 45939, $null, "09.10.202524"),
    @(160, 45940, $null, "10.10.20251"),
    @(161, 45940, $null, "10.10.20252"),
    @(162, 45940, $null, "10.10.20253"),
    @(163, 45940, $null, "10.10.20254"),
    @(164, 45940, $null, "10.10.20255"),
    @(165, 45940, $null, "10.10.20256"),
    @(166, 45940, $null, "10.10.20257"),
    @(167, 45940, 0.069, "10.10.20258"),
    @(168, 45940, 0.587, "10.10.20259"),
    @(169, 45940, 1.082, "10.10.202510"),
    @(170, 45940, 1.834, "10.10.202511")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
}
